$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: the stray "_GoBack" bookmark that used to sit at the end
# of the intro paragraph ("甲乙双方本着自愿平等...") is being moved
# elsewhere (see change 2). Word only ever keeps a single "_GoBack"
# bookmark, so re-adding it at the new location (below) automatically
# removes it from here - nothing else to do for this hunk on its own.
# ------------------------------------------------------------------

# ------------------------------------------------------------------
# Change 2: the table header cell that read "商标名称" (trademark
# name) becomes "专利名称" (patent name), typed as two runs - "专利"
# (newly inserted) followed by the relocated "_GoBack" bookmark and
# then the remaining "名称" text.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("商标", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($r.Find.Found) {
    # Replace "商标" with "专利" in place.
    $r.Text = "专利"

    # Force Word to keep this newly-typed text as its own run instead
    # of silently re-merging it with the following "名称" run (toggle
    # a format flag off/on so the run boundary sticks).
    $r.Font.Bold = $false
    $r.Font.Bold = $true

    # Re-seat the document's single "_GoBack" bookmark right after the
    # freshly typed "专利" text (i.e. between it and "名称"), exactly
    # like Word drops it at the last edit point.
    $bmRange = $d.Range($r.End, $r.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

"done"
